# Thu, Jun 04, 2020  4:04:57 AM
#
# 1. Re-style the table on slide 6 with a different built-in table style.
# 2. Re-colour the presentation's theme (previously the "Integral" palette)
#    to the standard "Office Theme" colour palette.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$tableShape = $p.Slides.Item(6).Shapes.Item(2)
$tableShape.Table.ApplyStyle("{12CB7BF4-9B34-4228-9FB1-392E0C1F2948}")

# --- 2. Theme colour scheme -------------------------------------------------
$tcs = $p.Designs.Item(1).SlideMaster.Theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0          # Dark 1    000000
$tcs.Colors(2).RGB  = 16777215   # Light 1   FFFFFF
$tcs.Colors(3).RGB  = 6968388    # Dark 2    44546A
$tcs.Colors(4).RGB  = 15132391   # Light 2   E7E6E6
$tcs.Colors(5).RGB  = 13998939   # Accent 1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # Accent 2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # Accent 3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # Accent 4  FFC000
$tcs.Colors(9).RGB  = 12874308   # Accent 5  4472C4
$tcs.Colors(10).RGB = 4697456    # Accent 6  70AD47
$tcs.Colors(11).RGB = 12673797   # Hyperlink 0563C1
$tcs.Colors(12).RGB = 7491477    # Followed Hyperlink 954F72
